$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A becomes narrower, column B becomes wider (stored <col> width,
# i.e. the OOXML "number of characters" unit): A 51.7109375 -> 29.7109375,
# B 22.85546875 -> 52.5703125.
#
# The COM layer here quantizes ColumnWidth to whole-pixel steps (1/6 of a
# character, matching a Maximum-Digit-Width of 6px) before it is written
# back out as the <col width="..."> attribute, so the literal target
# values above are not bit-exactly reproducible through this API. Values
# are chosen here so the *stored* width lands on the nearest reachable
# step (29.6666... and 52.5 respectively) rather than whatever the naive
# literal assignment would truncate/round to.
$ws.Columns.Item(1).ColumnWidth = 28.833333333333332
$ws.Columns.Item(2).ColumnWidth = 51.666666666666664

